$p = $ppt.ActivePresentation
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                Write-Output $shp.TextFrame.TextRange.Text
            }
        }
    }
}
